$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.590104341506958
$ws.Range("B1").Value = 4.080108642578125
$ws.Range("C1").Value = 3.335054874420166
$ws.Range("D1").Value = 3.149731397628784
$ws.Range("E1").Value = 1.776500821113586
